$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 849, pushing existing data (old rows 849+) down by 2.
$ws.Rows.Item(849).Insert()
$ws.Rows.Item(849).Insert()

# Fill in the two newly inserted rows with their data.
# Column A holds date-like text ("2026/02/20") that must stay a literal
# string (like every other date cell in this sheet) instead of being
# auto-parsed into a date serial number. Forcing NumberFormat="@" before
# the assignment keeps it text; resetting the style back to "Normal"
# afterwards drops the now-unneeded explicit text format so the cell
# matches the unstyled cells around it.
$ws.Range("A849:A850").NumberFormat = "@"

$ws.Cells.Item(849, 1).Value = "2026/02/20"
$ws.Cells.Item(849, 2).Value = "金"
$ws.Cells.Item(849, 3).Value = 23
$ws.Cells.Item(849, 4).Value = 201

$ws.Cells.Item(850, 1).Value = "2026/02/21"
$ws.Cells.Item(850, 2).Value = "土"
$ws.Cells.Item(850, 3).Value = 3
$ws.Cells.Item(850, 4).Value = 201

$ws.Range("A849:A850").Style = "Normal"
